$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Perejil (Feria Lagunitas de
# Puerto Montt). It belongs chronologically ahead of the existing rows, so
# insert a fresh row at 349 and push the existing data (rows 349-371) down to
# 350-372, then populate the new row with its values.
$ws.Rows("349:349").Insert()

$ws.Cells.Item(349, 1).Value = 4
$ws.Cells.Item(349, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(349, 3).Value = "Los Lagos"
$ws.Cells.Item(349, 4).Value = 45013
$ws.Cells.Item(349, 5).Value = 10
$ws.Cells.Item(349, 6).Value = 100112044
$ws.Cells.Item(349, 7).Value = "Perejil"
$ws.Cells.Item(349, 8).Value = "Sin especificar"
$ws.Cells.Item(349, 9).Value = "Primera"
$ws.Cells.Item(349, 10).Value = 160
$ws.Cells.Item(349, 11).Value = 6000
$ws.Cells.Item(349, 12).Value = 6000
$ws.Cells.Item(349, 13).Value = 6000
$ws.Cells.Item(349, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(349, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(349, 16).Value = 3000
$ws.Cells.Item(349, 17).Value = 2
$ws.Cells.Item(349, 18).Value = "Hortaliza"
